$wb = $excel.ActiveWorkbook
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    if ($lastRow -ge 2) {
        $colACells = $ws.Range($ws.Cells.Item(2,1), $ws.Cells.Item($lastRow,1))
        $colACells.Font.Bold = $true
    }
}
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count
    $hdr = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
    $hdr.Font.Color = 0
    $hdr.Font.Bold = $false
}
